$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the use-case step text (comment/documentation update): append the
# clarifying sentence about the employee acting on the information.
$ws.Range("G8").Value = "1. Customer makes a call and makes an order and specifies information (name,size,type) . Customer provides name and phone number. The Employee implements the information."

# The G8:G9 merged cell now needs more vertical room for the longer text -
# row 8 shrinks back to the default height and row 9 is manually expanded,
# matching how the author resized the row after editing the text.
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).RowHeight = 66

# Reflect the author's final cursor/selection position after the edit.
$null = $ws.Range("K9").Select()
